$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1285113333333333
$ws.Range("H2").Value = 0.385534
$ws.Range("I2").Value = 0.03749201237720504
$ws.Range("J2").Value = 0.03749201237720504
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 1.832234356757556
$ws.Range("R2").Value = 16.490109210818
$ws.Range("S2").Value = 0.007828147972399266
$ws.Range("T2").Value = 0.007828147972399269
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1285113333333333
$ws.Range("H3").Value = 0.385534
$ws.Range("I3").Value = 0.03749201237720504
$ws.Range("J3").Value = 0.03749201237720504
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 3.732326038809778
$ws.Range("R3").Value = 33.590934349288
$ws.Range("S3").Value = 0.01594621365180951
$ws.Range("T3").Value = 0.01594621365180951
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1285113333333333
$ws.Range("H4").Value = 0.385534
$ws.Range("I4").Value = 0.03749201237720504
$ws.Range("J4").Value = 0.03749201237720504
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 2.584471676758444
$ws.Range("R4").Value = 23.260245090826
$ws.Range("S4").Value = 0.01104205182133097
$ws.Range("T4").Value = 0.01104205182133097
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1285113333333333
$ws.Range("H5").Value = 0.385534
$ws.Range("I5").Value = 0.03749201237720504
$ws.Range("J5").Value = 0.03749201237720504
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 0.6262431809906667
$ws.Range("R5").Value = 5.636188628916
$ws.Range("S5").Value = 0.002675598931665289
$ws.Range("T5").Value = 0.00267559893166529
$ws.Range("I6").Value = 0.7552862722193517
$ws.Range("J6").Value = 0.755286272219352
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 36.91083431917934
$ws.Range("R6").Value = 332.197508872614
$ws.Range("S6").Value = 0.1577000626418678
$ws.Range("T6").Value = 0.1577000626418679
$ws.Range("I7").Value = 0.7552862722193517
$ws.Range("J7").Value = 0.755286272219352
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("Q7").Value = 75.18867198160267
$ws.Range("R7").Value = 676.698047834424
$ws.Range("S7").Value = 0.3212405923671146
$ws.Range("T7").Value = 0.3212405923671147
$ws.Range("I8").Value = 0.7552862722193517
$ws.Range("J8").Value = 0.755286272219352
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 52.06484940728865
$ws.Range("R8").Value = 468.5836446655979
$ws.Range("S8").Value = 0.2224449857179872
$ws.Range("T8").Value = 0.2224449857179873
$ws.Range("I9").Value = 0.7552862722193517
$ws.Range("J9").Value = 0.755286272219352
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 12.615830617852
$ws.Range("R9").Value = 113.542475560668
$ws.Range("S9").Value = 0.05390063149238206
$ws.Range("T9").Value = 0.05390063149238209
$ws.Range("G10").Value = 0.692415
$ws.Range("H10").Value = 2.077245
$ws.Range("I10").Value = 0.2020057770533527
$ws.Range("J10").Value = 0.2020057770533527
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 9.872020772235
$ws.Range("R10").Value = 88.84818695011501
$ws.Range("S10").Value = 0.04217781372051885
$ws.Range("T10").Value = 0.04217781372051887
$ws.Range("G11").Value = 0.692415
$ws.Range("H11").Value = 2.077245
$ws.Range("I11").Value = 0.2020057770533527
$ws.Range("J11").Value = 0.2020057770533527
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 20.10965466726
$ws.Range("R11").Value = 180.98689200534
$ws.Range("S11").Value = 0.08591769487815093
$ws.Range("T11").Value = 0.08591769487815096
$ws.Range("G12").Value = 0.692415
$ws.Range("H12").Value = 2.077245
$ws.Range("I12").Value = 0.2020057770533527
$ws.Range("J12").Value = 0.2020057770533527
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 13.925051663895
$ws.Range("R12").Value = 125.325464975055
$ws.Range("S12").Value = 0.05949422602312806
$ws.Range("T12").Value = 0.05949422602312808
$ws.Range("G13").Value = 0.692415
$ws.Range("H13").Value = 2.077245
$ws.Range("I13").Value = 0.2020057770533527
$ws.Range("J13").Value = 0.2020057770533527
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 3.37417845507
$ws.Range("R13").Value = 30.36760609563
$ws.Range("S13").Value = 0.01441604243155484
$ws.Range("T13").Value = 0.01441604243155484
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.01787866666666667
$ws.Range("H14").Value = 0.053636
$ws.Range("I14").Value = 0.005215938350090445
$ws.Range("J14").Value = 0.005215938350090446
$ws.Range("M14").Value = 14.25737566666667
$ws.Range("N14").Value = 42.772127
$ws.Range("O14").Value = 0.2087950866344732
$ws.Range("P14").Value = 0.2087950866344732
$ws.Range("Q14").Value = 0.2549028670857778
$ws.Range("R14").Value = 2.294125803772
$ws.Range("S14").Value = 0.001089062299687206
$ws.Range("T14").Value = 0.001089062299687206
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.01787866666666667
$ws.Range("H15").Value = 0.053636
$ws.Range("I15").Value = 0.005215938350090445
$ws.Range("J15").Value = 0.005215938350090446
$ws.Range("N15").Value = 87.128332
$ws.Range("O15").Value = 0.4253229592313036
$ws.Range("P15").Value = 0.4253229592313036
$ws.Range("Q15").Value = 0.5192461350168889
$ws.Range("R15").Value = 4.673215215152
$ws.Range("S15").Value = 0.002218458334228511
$ws.Range("T15").Value = 0.002218458334228512
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.01787866666666667
$ws.Range("H16").Value = 0.053636
$ws.Range("I16").Value = 0.005215938350090445
$ws.Range("J16").Value = 0.005215938350090446
$ws.Range("M16").Value = 20.11084633333333
$ws.Range("N16").Value = 60.332539
$ws.Range("O16").Value = 0.2945174484164121
$ws.Range("P16").Value = 0.2945174484164122
$ws.Range("Q16").Value = 0.3595551179782222
$ws.Range("R16").Value = 3.235996061804
$ws.Range("S16").Value = 0.001536184853965949
$ws.Range("T16").Value = 0.001536184853965949
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.01787866666666667
$ws.Range("H17").Value = 0.053636
$ws.Range("I17").Value = 0.005215938350090445
$ws.Range("J17").Value = 0.005215938350090446
$ws.Range("M17").Value = 4.873057999999999
$ws.Range("N17").Value = 14.619174
$ws.Range("O17").Value = 0.07136450571781097
$ws.Range("P17").Value = 0.07136450571781099
$ws.Range("Q17").Value = 0.08712377962933333
$ws.Range("R17").Value = 0.784114016664
$ws.Range("S17").Value = 0.0003722328622087791
$ws.Range("T17").Value = 0.0003722328622087793
